$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '35.852.24'
Set-TextValue $ws 'E2' '  -4.15%  '
Set-TextValue $ws 'D3' '1.962.87'
Set-TextValue $ws 'E3' '  -3.82%  '
Set-TextValue $ws 'E4' '  +0.17%  '
Set-TextValue $ws 'D5' '242.45'
Set-TextValue $ws 'E5' '  -4.03%  '
Set-TextValue $ws 'E6' '  -4.03%  '
Set-TextValue $ws 'D7' '61.83'
Set-TextValue $ws 'E7' '  -5.00%  '
Set-TextValue $ws 'E8' '  +0.06%  '
Set-TextValue $ws 'D9' '0.365'
Set-TextValue $ws 'E9' '  -2.73%  '
Set-TextValue $ws 'D10' '55.99'
Set-TextValue $ws 'E10' '  -5.24%  '
Set-TextValue $ws 'D11' '0.0797'
Set-TextValue $ws 'E11' '  +5.74%  '
Set-TextValue $ws 'E12' '  -1.44%  '
Set-TextValue $ws 'D13' '0.856'
Set-TextValue $ws 'E13' '  -5.55%  '
Set-TextValue $ws 'D14' '13.98'
Set-TextValue $ws 'E14' '  -7.68%  '
Set-TextValue $ws 'D15' '21.82'
Set-TextValue $ws 'E15' '  +5.39%  '
Set-TextValue $ws 'D16' '2.248.45'
Set-TextValue $ws 'E16' '  -3.78%  '
Set-TextValue $ws 'D17' '5.41'
Set-TextValue $ws 'E17' '  -3.53%  '
Set-TextValue $ws 'D18' '1.966.83'
Set-TextValue $ws 'E18' '  -3.48%  '
Set-TextValue $ws 'D19' '35.820.36'
Set-TextValue $ws 'E19' '  -3.94%  '
Set-TextValue $ws 'D20' '70.78'
Set-TextValue $ws 'E20' '  -3.23%  '
Set-TextValue $ws 'D21' '0.0₃0850'
Set-TextValue $ws 'E21' '  -2.62%  '
Set-TextValue $ws 'D22' '239.05'
Set-TextValue $ws 'E22' '  +0.95%  '
Set-TextValue $ws 'D23' '5.19'
Set-TextValue $ws 'E23' '  -3.02%  '
Set-TextValue $ws 'E24' '  -0.07%  '
Set-TextValue $ws 'E25' '  -8.80%  '
Set-TextValue $ws 'D26' '2.29'
Set-TextValue $ws 'E26' '  -2.68%  '
Set-TextValue $ws 'D27' '9.79'
Set-TextValue $ws 'E27' '  +2.30%  '
Set-TextValue $ws 'D28' '158.47'
Set-TextValue $ws 'E28' '  -4.49%  '
Set-TextValue $ws 'B29' 'Kaspa'
Set-TextValue $ws 'C29' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws 'D29' '0.132'
Set-TextValue $ws 'E29' '  +16.81%  '
Set-TextValue $ws 'B30' 'EthereumClassic'
Set-TextValue $ws 'C30' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws 'D30' '19.67'
Set-TextValue $ws 'E30' '  -1.05%  '
Set-TextValue $ws 'E31' '  -1.91%  '
Set-TextValue $ws 'D32' '4.86'
Set-TextValue $ws 'E32' '  -7.17%  '
Set-TextValue $ws 'D33' '1.14'
Set-TextValue $ws 'E33' '  -7.21%  '
Set-TextValue $ws 'D34' '0.0615'
Set-TextValue $ws 'E34' '  -0.10%  '
Set-TextValue $ws 'D35' '4.38'
Set-TextValue $ws 'E35' '  -7.35%  '
Set-TextValue $ws 'D36' '6.29'
Set-TextValue $ws 'E36' '  +5.44%  '
Set-TextValue $ws 'D37' '2.31'
Set-TextValue $ws 'E37' '  -6.43%  '
Set-TextValue $ws 'D38' '1.00'
Set-TextValue $ws 'E38' '  +0.16%  '
Set-TextValue $ws 'D39' '1.83'
Set-TextValue $ws 'E39' '  +1.22%  '
Set-TextValue $ws 'D40' '3.15'
Set-TextValue $ws 'E40' '  +15.21%  '
Set-TextValue $ws 'D41' '0.0981'
Set-TextValue $ws 'E41' '  -5.67%  '
Set-TextValue $ws 'D42' '1.23'
Set-TextValue $ws 'E42' '  -0.45%  '
Set-TextValue $ws 'E43' '  -3.54%  '
Set-TextValue $ws 'E44' '  -4.55%  '
Set-TextValue $ws 'D45' '1.09'
Set-TextValue $ws 'E45' '  -4.90%  '
Set-TextValue $ws 'B46' 'InjectiveProtocol'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws 'D46' '16.17'
Set-TextValue $ws 'E46' '  -5.66%  '
Set-TextValue $ws 'B47' 'Aave'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D47' '91.86'
Set-TextValue $ws 'E47' '  -3.72%  '
Set-TextValue $ws 'D48' '7.52'
Set-TextValue $ws 'E48' '  -7.91%  '
Set-TextValue $ws 'D49' '1.339.05'
Set-TextValue $ws 'E49' '  -5.80%  '
Set-TextValue $ws 'D50' '2.75'
Set-TextValue $ws 'E50' '  -6.07%  '
Set-TextValue $ws 'D51' '2.140.65'
Set-TextValue $ws 'E51' '  -3.86%  '
